$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (existing D:K data shifts to E:L)
$ws.Columns("D:D").Insert()

# Copy number formatting from the (now-shifted) E column into the new D column
# so the new column matches the formatting of the data it sits beside.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate new column D with the latest period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2188100
$ws.Range("D9").Value = 750600
$ws.Range("D10").Value = 1437600
$ws.Range("D12").Value = 113800
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 417500
$ws.Range("D17").Value = 1416300
$ws.Range("D18").Value = 771800
$ws.Range("D20").Value = -500
$ws.Range("D21").Value = 1188800
$ws.Range("D22").Value = 73200
$ws.Range("D23").Value = 698100
$ws.Range("D24").Value = 152500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 545700
$ws.Range("D27").Value = 545700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 11400
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 500
$ws.Range("D33").Value = 557000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 557000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 2300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 471700
$ws.Range("D44").Value = 11100
$ws.Range("D45").Value = 59500
$ws.Range("D46").Value = 544500
$ws.Range("D47").Value = 171300
$ws.Range("D48").Value = 3463600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 19300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4198800
$ws.Range("D57").Value = 30000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 257200
$ws.Range("D60").Value = 287300
$ws.Range("D61").Value = 1226100
$ws.Range("D62").Value = 597300
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2110700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1607700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2088200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 557000
$ws.Range("D83").Value = 417500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1104900
$ws.Range("D91").Value = -894500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -293400
$ws.Range("D96").Value = -111400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -1289300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -477800
